$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item Filter")

# Clear out the "Color" data (D2:D4) that's no longer needed by this test step.
$ws.Range("D2:D4").ClearContents() | Out-Null

# Reflect the new selection left behind by the edit.
$ws.Range("D2:D4").Select() | Out-Null
